# Kreiranje korisničkog računa - use-case scenario sheet
# Dodan scenarij upotrebe aplikacije od strane uposlenika
#
# 1) Text corrections (double spaces -> single/trailing spaces, hyphen -> en-dash)
# 2) Visual overhaul: bold labels in column A, bordered table (medium box grid),
#    wrap text + vertical-center alignment throughout the table
# 3) Two blank spacer rows appended under the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the text content -------------------------------------------------
$ws.Range("B2").Value = "Korisnik putem aplikacije podnosi zahtjev za otvaranje korisničkog računa "
$ws.Range("A4").Value = "Posljedice " + [char]0x2013 + " uspješan završetak"
$ws.Range("A5").Value = "Posljedice " + [char]0x2013 + " neuspješan završetak"
$ws.Range("B6").Value = "Korisnik, administrator "
$ws.Range("B7").Value = "Korisnik podnosi zahtjev za formiranje korisničkog računa, popuni odgovarajuće podatke, otvori korisnički račun"
